$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("anemia prevalence")

# Remove the old placeholder comment on L2 (no longer applicable once the
# sheet is restructured as a distribution rather than a single prevalence
# row)
$ws.Range("L2").Comment.Delete()

# Shift the existing prevalence table two columns to the right, to make room
# for new "Distribution" / "Status" label columns (matching the layout used
# on the "distributions" sheet)
$ws.Range("A1:B1").EntireColumn.Insert()

# New label columns, matching the header style used for the rest of row 1
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Distribution"
$ws.Range("B1").Value = "Status"

$ws.Range("A2").Value = "Anemia"
$ws.Range("B2").Value = "anemic"

# The data used to be entered as a prevalence fraction (0.5 == 50%); now it
# is entered as a percentage distribution value (50)
$ws.Range("C2:N2").Value = 50

# Add the complementary "not anemic" distribution row
$ws.Range("B3").Value = "not anemic"
$ws.Range("C3").Formula = "=100-C2"
$ws.Range("D3:N3").Formula = "=100-D2"

$ws.Range("H9").Select()

$wb.Worksheets.Item("distributions").Range("D12").Select()
$ws.Select()
